# Daily 04.07 took place
# Fill in the "04.07" row (row 3, date 45111) for every team member's sheet
# with their yesterday/today/blocker entries, then leave the "Raik" sheet
# active (mirrors where the author ended up working last).

$wb = $excel.ActiveWorkbook

$tom = $wb.Worksheets.Item("Tom")
$tom.Range("B3").Value = "starting with the website "
$tom.Range("C3").Value = "working on the website "
$tom.Range("D3").Value = "limited knwoledge"

$linus = $wb.Worksheets.Item("Linus")
$linus.Range("B3").Value = "scrum organisation, helping Raik with the database"
$linus.Range("C3").Value = "database, hopefully starting with the hardware"
$linus.Range("D3").Value = "sql syntax"

$raik = $wb.Worksheets.Item("Raik")
$raik.Range("B3").Value = "database, erm and dump"
$raik.Range("C3").Value = "working on the sql dump, hardware"
$raik.Range("D3").Value = "sql syntax of course"

$arweed = $wb.Worksheets.Item("Arweed")
$arweed.Range("B3").Value = "starting with the website"
$arweed.Range("C3").Value = "VM problems solving, Linux things"
$arweed.Range("D3").Value = "Linux"

# The now-unused placeholder cells (rows 4-9) on Linus/Raik/Arweed lose their
# leftover style formatting once the sheet is tidied up after entry.
$linus.Range("B4:D9").Clear()
$raik.Range("B4:D9").Clear()
$arweed.Range("B4:D9").Clear()

# Columns auto-widen to fit the freshly typed text.
$linus.Columns("B").ColumnWidth = 40.1666666666667
$linus.Columns("C").ColumnWidth = 37.1666666666667
$linus.Columns("D").ColumnWidth = 22.0

$raik.Columns("B").ColumnWidth = 21.5
$raik.Columns("C").ColumnWidth = 28.1666666666667
$raik.Columns("D").ColumnWidth = 22.0

$arweed.Columns("B").ColumnWidth = 21.5
$arweed.Columns("C").ColumnWidth = 26.3333333333333
$arweed.Columns("D").ColumnWidth = 22.0

# Per-sheet selection state, matching where each user last clicked.
$tom.Range("D3").Select() | Out-Null
$linus.Range("D4").Select() | Out-Null
$raik.Range("D4").Select() | Out-Null
$arweed.Range("B3").Select() | Out-Null

# Raik ends up as the active sheet/tab.
$raik.Activate() | Out-Null
